# Update "想去人数" (interested count) values in column F
# for both the "展览" and "全部类型" sheets, which hold identical data.

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F
$updates = @{
    3  = 48
    5  = 85
    6  = 127
    7  = 1241
    8  = 1523
    10 = 385
    12 = 143
    17 = 297
    19 = 1718
    20 = 67
    23 = 661
    26 = 4141
    28 = 262
    29 = 1080
    32 = 505
    34 = 229
    36 = 135
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
